$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column I header
$ws.Range("I1").Value = 'Other found locations'

# Row 2
$ws.Range("E2").Value = '[Na%Zhu%NULL%0,   Dingyu%Zhang%NULL%0,   Wenling%Wang%NULL%0,   Xingwang%Li%NULL%0,   Bo%Yang%NULL%0,   Jingdong%Song%NULL%0,   Xiang%Zhao%NULL%0,   Baoying%Huang%NULL%0,   Weifeng%Shi%NULL%0,   Roujian%Lu%NULL%0,   Peihua%Niu%NULL%0,   Faxian%Zhan%NULL%0,   Xuejun%Ma%NULL%0,   Dayan%Wang%NULL%0,   Wenbo%Xu%NULL%0,   Guizhen%Wu%NULL%0,   George F.%Gao%NULL%0,   Wenjie%Tan%NULL%0]'
$ws.Range("I2").Value = '_PMC'

# Row 3
$ws.Range("E3").Value = '[Xiaoxia%Lu%NULL%1,   Liqiong%Zhang%NULL%1,   Hui%Du%NULL%0,   Jingjing%Zhang%NULL%2,   Jingjing%Zhang%NULL%0,   Yuan Y.%Li%NULL%1,   Jingyu%Qu%NULL%1,   Wenxin%Zhang%NULL%1,   Youjie%Wang%NULL%1,   Shuangshuang%Bao%NULL%1,   Ying%Li%NULL%1,   Chuansha%Wu%NULL%1,   Hongxiu%Liu%NULL%1,   Di%Liu%NULL%2,   Di%Liu%NULL%0,   Jianbo%Shao%NULL%3,   Jianbo%Shao%NULL%0,   Jianbo%Shao%NULL%0,   Xuehua%Peng%NULL%1,   Yonghong%Yang%NULL%2,   Yonghong%Yang%NULL%0,   Zhisheng%Liu%NULL%3,   Zhisheng%Liu%NULL%0,   Zhisheng%Liu%NULL%0,   Yun%Xiang%NULL%1,   Furong%Zhang%NULL%1,   Rona M.%Silva%NULL%2,   Rona M.%Silva%NULL%0,   Kent E.%Pinkerton%NULL%1,   Kunling%Shen%NULL%2,   Kunling%Shen%NULL%0,   Han%Xiao%NULL%0,   Han%Xiao%NULL%0,   Shunqing%Xu%NULL%2,   Shunqing%Xu%NULL%0,   Gary W.K.%Wong%NULL%2,   Gary W.K.%Wong%NULL%0]'
$ws.Range("I3").Value = '_PMC'

# Row 4
$ws.Range("I4").Value = ""

# Row 5
$ws.Range("E5").Value = '[Esquivel%Sosa Leidelen%coreGivesNoEmail%1,  Mart\u00ednez-Fort\u00fan%Amador Maryla%coreGivesNoEmail%1,  \u00c1guila%Carbelo Madyaret%coreGivesNoEmail%1]'
$ws.Range("F5").Value = "not found"
$ws.Range("G5").Value = "N/A"
$ws.Range("I5").Value = ""

# Row 6
$ws.Range("I6").Value = '_PMC'

# Row 7
$ws.Range("I7").Value = ""

# Row 8
$ws.Range("E8").Value = '[Ruiyun%Li%NULL%1,   Sen%Pei%NULL%2,   Sen%Pei%NULL%0,   Bin%Chen%NULL%2,   Bin%Chen%NULL%0,   Yimeng%Song%NULL%2,   Yimeng%Song%NULL%0,   Tao%Zhang%NULL%2,   Tao%Zhang%NULL%0,   Wan%Yang%NULL%1,   Jeffrey%Shaman%NULL%2,   Jeffrey%Shaman%NULL%0]'
$ws.Range("I8").Value = '_PMC'

# Row 9
$ws.Range("E9").Value = '[Andrea%T. Cruz%coreGivesNoEmail%1,  Steven%L. Zeichner%coreGivesNoEmail%1]'
$ws.Range("F9").Value = "not found"
$ws.Range("G9").Value = "N/A"
$ws.Range("I9").Value = ""

# Row 10
$ws.Range("E10").Value = '[Dan%Sun%NULL%2,   Hui%Li%NULL%0,   Xiao-Xia%Lu%NULL%2,   Han%Xiao%NULL%0,   Jie%Ren%NULL%2,   Fu-Rong%Zhang%792523496@qq.com%2,   Zhi-Sheng%Liu%liuzsc@126.com%2]'
$ws.Range("I10").Value = '_PMC_Springer'

# Row 11
$ws.Range("E11").Value = '[Hin%Chu%kyyuen@hku.hk%1,   Jie%Zhou%NULL%1,   Bosco Ho-Yin%Wong%NULL%1,   Cun%Li%NULL%1,   Jasper Fuk-Woo%Chan%NULL%1,   Zhong-Shan%Cheng%NULL%1,   Dong%Yang%NULL%1,   Dong%Wang%NULL%1,   Andrew Chak-Yiu%Lee%NULL%1,   Chuangen%Li%NULL%1,   Man-Lung%Yeung%NULL%1,   Jian-Piao%Cai%NULL%0,   Ivy Hau-Yee%Chan%NULL%1,   Wai-Kuen%Ho%NULL%1,   Kelvin Kai-Wang%To%NULL%0,   Bo-Jian%Zheng%NULL%1,   Yanfeng%Yao%NULL%1,   Chuan%Qin%NULL%0,   Kwok-Yung%Yuen%NULL%0]'
$ws.Range("I11").Value = '_PMC'

# Row 12
$ws.Range("E12").Value = '[William J.%Liu%liujun@ivdc.chinacdc.cn%1,   Min%Zhao%NULL%0,   Kefang%Liu%NULL%1,   Kun%Xu%NULL%1,   Gary%Wong%NULL%1,   Wenjie%Tan%NULL%0,   George F.%Gao%gaofu@chinacdc.cn%0]'
$ws.Range("I12").Value = '_PMC_elsevier'

# Row 13
$ws.Range("E13").Value = '[Henry%Brandon Michael%coreGivesNoEmail%1,  Lippi%Giuseppe%coreGivesNoEmail%6,  Plebani%Mario%coreGivesNoEmail%1]'
$ws.Range("F13").Value = "not found"
$ws.Range("G13").Value = "N/A"
$ws.Range("I13").Value = ""
